# "updated file at 20.05.2017 in 20 min"
#
# Row 62: the previously "IN PROGRESS" Sell order is finalized -> DONE, with
#         Finalized date / Fee / Profit(%) / Transaction duration filled in,
#         and its Current value / Transaction value refreshed.
# Row 63: a brand new "IN PROGRESS" Buy order is logged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Row 62 - finalize the existing Sell order
# ---------------------------------------------------------------------

# D62 ("Current value") keeps its text but must stay TEXT (it contains
# embedded newlines that Excel would otherwise read as numeric) -
# force text format first, write it, then restore the normal wrapped
# "current value" look (style copied from the cell above it).
$ws.Range("D62").NumberFormat = "@"
$ws.Range("D62").Value = "              0.337`n`n"
$ws.Range("D61").Copy()
$ws.Range("D62").PasteSpecial(-4122)
$ws.Rows.Item(62).RowHeight = 14.25

$ws.Range("E62").Value = "         0.350  USDT"
$ws.Range("H62").Value = "DONE"

$ws.Range("I62").Value = 42875.454837962963

$ws.Range("J62").Value = "0.06342971 USDT (0.15%)"
$ws.Range("L62").Value = "2.5 day"

# K62 ("Profit(%)") is rendered as "     " + a green "~4%"
$ws.Range("K62").Value = "     ~4%"
$k62 = $ws.Range("K62").Characters(6, 3)
$k62.Font.Color = 5287936
$k62.Font.Name = "Calibri"
$k62.Font.Size = 11

# ---------------------------------------------------------------------
# Row 63 - log the new Buy order
# ---------------------------------------------------------------------

$ws.Range("A63").Value = 42875.451099537036
$ws.Range("A62").Copy()
$ws.Range("A63").PasteSpecial(-4122)

# B63 ("Action") = green "Buy", same look used throughout column B
$ws.Range("B63").Value = "            Buy"
$b63 = $ws.Range("B63").Characters(13, 3)
$b63.Font.Color = 5287936
$b63.Font.Name = "Calibri"
$b63.Font.Size = 11

$ws.Range("C63").Value = "        XRP"

$ws.Range("D63").NumberFormat = "@"
$ws.Range("D63").Value = "              0.355`n"
$ws.Range("D61").Copy()
$ws.Range("D63").PasteSpecial(-4122)

$ws.Range("E63").Value = "         0.300  USDT"
$ws.Range("F63").Value = "         145 XRP"
$ws.Range("G63").Value = " XRP/USDT0000009"
$ws.Range("H63").Value = "IN PROGRESS"

$ws.Range("I61").Copy()
$ws.Range("I63").PasteSpecial(-4122)

$ws.Range("K63").Value = "     "

$ws.Rows.Item(63).RowHeight = 14.25

# ---------------------------------------------------------------------
# Selection / view bookkeeping to match the author's saved session
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 44
$ws.Range("F63").Select() | Out-Null
